# Update the ANOVA table with the results from the 202051023 classification
# model run. Each replacement targets a unique whole-word match so that the
# correct cell is updated even though the document contains repeated tokens
# (e.g. "2", "NA") elsewhere in the table.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "0.011"; New = "0.001" },
    @{ Old = "1.656"; New = "0.428" },
    @{ Old = "0.224"; New = "0.656" },
    @{ Old = "0.005"; New = "0.001" },
    @{ Old = "1.554"; New = "0.353" },
    @{ Old = "0.232"; New = "0.558" },
    @{ Old = "0.002"; New = "0.000" },
    @{ Old = "0.248"; New = "0.093" },
    @{ Old = "0.784"; New = "0.912" },
    @{ Old = "0.049"; New = "0.035" },
    @{ Old = "15";    New = "24" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $r.New, 2)
}
